$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 104, shifting existing rows 104:209 down to 105:210
$ws.Rows("104:104").Insert()

# Populate the newly inserted row 104 with the new record
$ws.Range("A104").Value = 4
$ws.Range("B104").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C104").Value = "Los Lagos"
$ws.Range("D104").Value = 44586
$ws.Range("E104").Value = 10
$ws.Range("F104").Value = "Fruta"
$ws.Range("G104").Value = 100104
$ws.Range("H104").Value = "Frutos de pepita"
$ws.Range("I104").Value = 100104005
$ws.Range("J104").Value = "Pera"
$ws.Range("K104").Value = "Packham's Triumph"
$ws.Range("L104").Value = "Primera"
$ws.Range("M104").Value = 500
$ws.Range("N104").Value = 14000
$ws.Range("O104").Value = 14000
$ws.Range("P104").Value = 14000
$ws.Range("Q104").Value = "$/caja 15 kilos empedrada"
$ws.Range("R104").Value = "Región de O'Higgins"
$ws.Range("S104").Value = 933
$ws.Range("T104").Value = 15
